$wb = $excel.ActiveWorkbook

# Copy the "Croatia" sheet (last sheet) to create the new "Greece" sheet, placed right after it.
$croatia = $wb.Worksheets.Item("Croatia")
$croatia.Copy([System.Reflection.Missing]::Value, $croatia)

# The newly created copy becomes the active sheet, placed right after Croatia.
$greece = $wb.ActiveSheet
$greece.Name = "Greece"

# Update the two changed cells on the new sheet (NGC code first, then market name,
# so new shared-string entries land in the same order as the target workbook).
$greece.Range("B4").Value = "NGC-4119/T3167/T3166"
$greece.Range("B2").Value = "Greece Market"
[void]$greece.Range("D16").Select()

# The old Croatia sheet loses its tab-selected state and its selection becomes the
# entire sheet (all columns selected).
[void]$croatia.Cells.Select()

# Re-activate Greece so it ends up as the active/selected tab.
[void]$greece.Activate()
